# The source feed published a new quarterly data point, so a new row is
# inserted at the top of the data table (row 2), pushing every existing
# row down by one. Excel automatically re-targets relative formula
# references and shared-formula ranges when rows are inserted this way.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the latest observation.
$ws.Range("A2").Value = 45657
$ws.Range("B2").Value = 72.34
$ws.Range("C2").Formula = "=(B2/B14-1)*100"

# The blank row Excel inserts copies formatting from the row above (the
# bold header row). Re-apply the formatting used by the rest of the data
# rows by copying it down from the row immediately below.
$ws.Range("A3:C3").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
